# ValueSet-cibmtr-priority-variables-epic-2021.xlsx
# Apply the "Metadata" sheet updates described by the commit:
#   - Version 0.1.6 -> 0.1.7
#   - Status active -> draft
#   - Date refreshed
#   - Contact split into org contact (w/ URL) + a named contact (Bob Milius)
#   - New "Jurisdiction" row inserted (blank value)
#   - Remaining rows (Description/Purpose/Copyright/Immutable) shift down one row

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- Extend the sheet by one row, carrying the existing formatting (style "2")
# of row 15 down into the brand new row 16, so every shifted row keeps the
# right look (vertical-top, wrap, bordered).
$ws.Range("A15:B15").Copy($ws.Range("A16:B16"))

# --- Shift rows 11..15 down into 12..16 (walk bottom-up so we never
# clobber a source row before it has been read).
for ($r = 15; $r -ge 11; $r--) {
    $dst = $r + 1
    $ws.Range("A$dst").Value = $ws.Range("A$r").Value2
    $bVal = $ws.Range("B$r").Value2
    if ($null -eq $bVal) {
        $ws.Range("B$dst").ClearContents()
    } else {
        $ws.Range("B$dst").Value = $bVal
    }
}

# --- Now overwrite the "top" rows with their new content.

# Row 3: Version
$ws.Range("B3").Value = "0.1.7"

# Row 6: Status
$ws.Range("B6").Value = "draft"

# Row 8: Date
$ws.Range("B8").Value = "2024-08-27T12:23:18-05:00"

# Row 10: Contact (org, now with URL)
$ws.Range("A10").Value = "Contact"
$ws.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"

# Row 11: Contact (named person) -- replaces the shifted-down duplicate
$ws.Range("A11").Value = "Contact"
$ws.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"

# Row 12: new Jurisdiction row (blank value) -- replaces the shifted-down Description
$ws.Range("A12").Value = "Jurisdiction"
$ws.Range("B12").ClearContents()

# Rows 13-16 already hold the correct shifted-down content:
#   13 Description | Priority Variables for Epic CIBMTR Reporting App (FY21)
#   14 Purpose      | (blank)
#   15 Copyright    | (blank)
#   16 Immutable    | BooleanType[null]

Write-Output "Metadata sheet updated"
